$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 68: new entry for "Aula 61 - Convertendo String para Integer e validando" ---

# Column B (aula number) - copy formatting from the row above (style index 5) and set value
$ws.Range("B65").Copy()
$ws.Range("B68").PasteSpecial(-4122)
$ws.Range("B68").Value = 61

# Column C (sessão) - same formatting/style, reuse existing session label text
$ws.Range("C65").Copy()
$ws.Range("C68").PasteSpecial(-4122)
$ws.Range("C68").Value = $ws.Range("C65").Value2

# Column D (nome da aula) - wrap-text formatting (style index 1)
$ws.Range("D65").Copy()
$ws.Range("D68").PasteSpecial(-4122)
$ws.Range("D68").Value = "61. Convertendo String para Integer e Validando"

# Column E (observação) - wrap-text formatting (style index 1)
$ws.Range("E65").Copy()
$ws.Range("E68").PasteSpecial(-4122)
$ws.Range("E68").Value = "criado uma classe StringToInteger para validar o campo de ""numero de endereço"" proibindo salvar ou submeter o formulario caso o usuario digite letras no campo de numero"

# --- Row 70: stray empty underlined cell left by the author at E70 ---
$ws.Range("E70").Font.Underline = $true

# --- Selection as left by the author after editing ---
$ws.Range("E65").Select() | Out-Null
